# Generate Report for Handback
# Updates the localization-status report after a handback event:
#  - "Ready for handoff" -> "Handed back: in sync with en-US" (status cells)
#  - populates "Latest Target File" / "Latest Handback File" / "Latest
#    Handback DateTime" for row 2/3 on the zh-cn and de-de detail sheets
#  - widens a couple of columns that now hold longer text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns (zh-cn / de-de) for both rows ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Overview: widen the zh-cn / de-de status columns to fit the longer text
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn detail sheet ---
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 12:37:08"
$wsZh.Range("K3").Value = "2016-08-19 12:37:08"

$zhALink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77fa9888f2ef475e4b96c658fee4001221fb6ac3/e2e/a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhALink, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhALink, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(10).ColumnWidth = 40

# --- de-de detail sheet ---
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 12:37:14"
$wsDe.Range("K3").Value = "2016-08-19 12:37:14"

$deALink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77fa9888f2ef475e4b96c658fee4001221fb6ac3/e2e/a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deALink, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deALink, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(10).ColumnWidth = 40
